$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row heights for rows 1 and 2 (new customHeight)
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(2).RowHeight = 15

# Copy formatting from column Q into column R for rows 3-38, then set values
# where the diff specifies a value (numeric rows); leave formatting-only cells blank.
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null

$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null

$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 88.796593100633856

$ws.Range("Q7").Copy() | Out-Null
$ws.Range("R7").PasteSpecial(-4122) | Out-Null
$ws.Range("R7").Value = 86.908583391486388

$ws.Range("Q8").Copy() | Out-Null
$ws.Range("R8").PasteSpecial(-4122) | Out-Null
$ws.Range("R8").Value = 89.680106631122953

$ws.Range("Q9").Copy() | Out-Null
$ws.Range("R9").PasteSpecial(-4122) | Out-Null
$ws.Range("R9").Value = 95.775910364145659

$ws.Range("Q10").Copy() | Out-Null
$ws.Range("R10").PasteSpecial(-4122) | Out-Null
$ws.Range("R10").Value = 96.517042279754136

$ws.Range("Q11").Copy() | Out-Null
$ws.Range("R11").PasteSpecial(-4122) | Out-Null
$ws.Range("R11").Value = 90.311530128242666

$ws.Range("Q12").Copy() | Out-Null
$ws.Range("R12").PasteSpecial(-4122) | Out-Null
$ws.Range("R12").Value = 90.746324915190343

$ws.Range("Q13").Copy() | Out-Null
$ws.Range("R13").PasteSpecial(-4122) | Out-Null
$ws.Range("R13").Value = 90.894107952204379

$ws.Range("Q14").Copy() | Out-Null
$ws.Range("R14").PasteSpecial(-4122) | Out-Null
$ws.Range("R14").Value = 81.065680730752504

$ws.Range("Q15").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Value = 85.088888888888889

$ws.Range("Q16").Copy() | Out-Null
$ws.Range("R16").PasteSpecial(-4122) | Out-Null

$ws.Range("Q17").Copy() | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").Value = 93.37839883628321

$ws.Range("Q18").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null
$ws.Range("R18").Value = 93.091416608513612

$ws.Range("Q19").Copy() | Out-Null
$ws.Range("R19").PasteSpecial(-4122) | Out-Null
$ws.Range("R19").Value = 94.815061646117954

$ws.Range("Q20").Copy() | Out-Null
$ws.Range("R20").PasteSpecial(-4122) | Out-Null
$ws.Range("R20").Value = 100.53781512605042

$ws.Range("Q21").Copy() | Out-Null
$ws.Range("R21").PasteSpecial(-4122) | Out-Null
$ws.Range("R21").Value = 100.33525796237662

$ws.Range("Q22").Copy() | Out-Null
$ws.Range("R22").PasteSpecial(-4122) | Out-Null
$ws.Range("R22").Value = 93.78989283832054

$ws.Range("Q23").Copy() | Out-Null
$ws.Range("R23").PasteSpecial(-4122) | Out-Null
$ws.Range("R23").Value = 95.401432340746325

$ws.Range("Q24").Copy() | Out-Null
$ws.Range("R24").PasteSpecial(-4122) | Out-Null
$ws.Range("R24").Value = 92.308748798242007

$ws.Range("Q25").Copy() | Out-Null
$ws.Range("R25").PasteSpecial(-4122) | Out-Null
$ws.Range("R25").Value = 89.338842975206617

$ws.Range("Q26").Copy() | Out-Null
$ws.Range("R26").PasteSpecial(-4122) | Out-Null
$ws.Range("R26").Value = 87.955555555555549

$ws.Range("Q27").Copy() | Out-Null
$ws.Range("R27").PasteSpecial(-4122) | Out-Null

$ws.Range("Q28").Copy() | Out-Null
$ws.Range("R28").PasteSpecial(-4122) | Out-Null
$ws.Range("R28").Value = 89.631204460036727

$ws.Range("Q29").Copy() | Out-Null
$ws.Range("R29").PasteSpecial(-4122) | Out-Null
$ws.Range("R29").Value = 89.204466154919743

$ws.Range("Q30").Copy() | Out-Null
$ws.Range("R30").PasteSpecial(-4122) | Out-Null
$ws.Range("R30").Value = 84.751749416861045

$ws.Range("Q31").Copy() | Out-Null
$ws.Range("R31").PasteSpecial(-4122) | Out-Null
$ws.Range("R31").Value = 96.201680672268907

$ws.Range("Q32").Copy() | Out-Null
$ws.Range("R32").PasteSpecial(-4122) | Out-Null
$ws.Range("R32").Value = 95.567144719687093

$ws.Range("Q33").Copy() | Out-Null
$ws.Range("R33").PasteSpecial(-4122) | Out-Null
$ws.Range("R33").Value = 91.330444457457389

$ws.Range("Q34").Copy() | Out-Null
$ws.Range("R34").PasteSpecial(-4122) | Out-Null
$ws.Range("R34").Value = 91.368262344515642

$ws.Range("Q35").Copy() | Out-Null
$ws.Range("R35").PasteSpecial(-4122) | Out-Null
$ws.Range("R35").Value = 92.345373803964662

$ws.Range("Q36").Copy() | Out-Null
$ws.Range("R36").PasteSpecial(-4122) | Out-Null
$ws.Range("R36").Value = 88.660287081339717

$ws.Range("Q37").Copy() | Out-Null
$ws.Range("R37").PasteSpecial(-4122) | Out-Null
$ws.Range("R37").Value = 84.944444444444443

$ws.Range("Q38").Copy() | Out-Null
$ws.Range("R38").PasteSpecial(-4122) | Out-Null

# Clear clipboard marching ants / copy mode
$excel.CutCopyMode = 0

# Update the active selection to R3, matching the authored selection change
$ws.Range("R3").Select() | Out-Null
